# Update the dated heading (2023-12-30 Saturday -> 2023-12-31 Sunday)
$d = $word.ActiveDocument

$p1 = $d.Paragraphs.Item(1)
$null = $p1.Range.Find.Execute("2023-12-30 Saturday", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "2023-12-31 Sunday", 2)

# Update the two-digit division problems in the single worksheet table.
# Cells are addressed by explicit (row, column) position so that values
# which coincidentally repeat elsewhere in the grid (e.g. "90÷4=" and
# "94÷6=" appear both as an original value in one cell and as a new
# value in another) are never double-replaced.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; New = "92÷8=" },
    @{ Row = 1;  Col = 2; New = "55÷2=" },
    @{ Row = 1;  Col = 3; New = "48÷5=" },
    @{ Row = 1;  Col = 4; New = "70÷4=" },
    @{ Row = 1;  Col = 5; New = "90÷4=" },

    @{ Row = 5;  Col = 1; New = "58÷3=" },
    @{ Row = 5;  Col = 2; New = "35÷8=" },
    @{ Row = 5;  Col = 3; New = "20÷7=" },
    @{ Row = 5;  Col = 4; New = "63÷3=" },
    @{ Row = 5;  Col = 5; New = "38÷4=" },

    @{ Row = 9;  Col = 1; New = "55÷3=" },
    @{ Row = 9;  Col = 2; New = "12÷2=" },
    @{ Row = 9;  Col = 3; New = "47÷2=" },
    @{ Row = 9;  Col = 4; New = "50÷2=" },
    @{ Row = 9;  Col = 5; New = "55÷8=" },

    @{ Row = 13; Col = 1; New = "94÷6=" },
    @{ Row = 13; Col = 2; New = "41÷5=" },
    @{ Row = 13; Col = 3; New = "41÷4=" },
    @{ Row = 13; Col = 4; New = "16÷7=" },
    @{ Row = 13; Col = 5; New = "35÷2=" },

    @{ Row = 17; Col = 1; New = "88÷7=" },
    @{ Row = 17; Col = 2; New = "72÷4=" },
    @{ Row = 17; Col = 3; New = "72÷2=" },
    @{ Row = 17; Col = 4; New = "79÷2=" },
    @{ Row = 17; Col = 5; New = "61÷7=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}
